$d = $word.ActiveDocument
$bullet = [char]0x2022

# --- Change 1: Collapse the three CORE COMPETENCIES detail paragraphs into one ---
$coreP = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Data Visualization & Design: Interactive Dashboards")) {
        $coreP = $p
        break
    }
}

$coreIdx = $coreP.Index
$coreP.Range.Text = "Data Visualization & Design " + $bullet + " Geospatial Analysis & Mapping " + $bullet + " Technical Visualization"

# The two following paragraphs (previously "Geospatial Analysis & Mapping: ..." and
# "Technical Visualization: ...") are now redundant - remove them.
$d.Paragraphs.Item($coreIdx + 1).Range.Delete()
$d.Paragraphs.Item($coreIdx + 1).Range.Delete()

# --- Change 2: Insert a new TECHNICAL SKILLS section before the closing line ---
$ledP = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("$bullet Led multi-million dollar research projects")) {
        $ledP = $p
        break
    }
}

$ledIdx = $ledP.Index
$ledP.Range.InsertParagraphAfter()

$headingP = $d.Paragraphs.Item($ledIdx + 1)
$headingP.Range.Text = "TECHNICAL SKILLS"
$headingP.Range.Style = "Heading2"

$headingP.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Item($ledIdx + 2)
$p1.Range.Style = "Normal"
$p1.Range.Text = "DATA VISUALIZATION & DESIGN Interactive Dashboards; Statistical Visualization; Geospatial Mapping; Choropleth Design"

$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($ledIdx + 3)
$p2.Range.Style = "Normal"
$p2.Range.Text = "GEOSPATIAL ANALYSIS & MAPPING Spatial Analysis; Mapping Technologies; Web Mapping; Spatial Data Processing"

$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item($ledIdx + 4)
$p3.Range.Style = "Normal"
$p3.Range.Text = "TECHNICAL VISUALIZATION Programming; Database Integration; Web Technologies; Statistical Computing"

Write-Output "Done"
